# Adding discussed definitions to correct files
#
# Marks specific rows with an asterisk ("*") in column A to flag the
# entries whose definitions were discussed, and restores the final
# selection to the bottom-right pane's last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToFlag = @(3, 4, 5, 6, 7, 8, 16, 17, 18, 19, 29, 31, 32)

foreach ($r in $rowsToFlag) {
    $ws.Cells.Item($r, 1).Value = "*"
}

$ws.Range("E33").Select()
